$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new values are plain numeric strings (e.g. "97.40").
# Excel would normally auto-convert these to numbers on assignment, which
# would lose the original text formatting (trailing zeros, precision, etc).
# Force text format first, assign, then restore default styling so the
# cell keeps looking like the untouched neighboring cells (no explicit
# number format / style index lingering on it).
$priceCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D15", "D17", "D19", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D41", "D42", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Assign every changed cell its new value, in the same order as the source diff.
$ws.Range("D2").Value = "42.673.24"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.544.41"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "308.93"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").Value = "97.40"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "35.44"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "0.0805"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "7.39"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").Value = "2.934.99"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "15.70"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").Value = "2.560.60"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "0.833"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").Value = "42.703.95"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "12.36"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "69.18"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "247.01"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("D24").Value = "2.90"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "26.52"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "40.21"
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("D31").Value = "157.40"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "5.71"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").Value = "0.0793"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "3.27"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").Value = "2.07"
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("D37").Value = "18.38"
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").Value = "2.58"
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "22.20"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").Value = "4.03"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "0.0298"
$ws.Range("D45").Value = "1.989.86"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "3.19"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("D47").Value = "9.04"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "2.789.40"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "80.63"
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.192"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "73.10"
$ws.Range("E51").Value = "  -3.60%  "

# Restore default (un-styled) formatting on the price cells we forced to text,
# so they match the plain/no-style cells surrounding them in the sheet.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
